$d = $word.ActiveDocument

function Replace-ParaText($index, $oldText, $newText) {
    $p = $d.Paragraphs($index)
    $r = $p.Range
    $found = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Host "WARNING: replace not found at paragraph $index : $oldText"
    }
}

function Insert-ParaAfter($index, $text) {
    $d.Paragraphs($index).Range.InsertParagraphAfter()
    $d.Paragraphs($index + 1).Range.Text = $text
}

# ---------------------------------------------------------------------------
# Work from the BOTTOM of the document UPWARD so that original (1-based)
# paragraph indices above the current edit point stay valid.
# ---------------------------------------------------------------------------

# Paragraph 44: df.select("Source", "Destination", "classification_final")...
Replace-ParaText 44 `
    '("Source", "Destination", "classification_final").show(truncate=False)' `
    '("cleaned_name", "cleaned_name_modified", "classification_final").show(truncate=False)'

# Paragraph 40: remove the blank paragraph that sat between ")" and
# "# Drop temporary row_number column"
$d.Paragraphs(40).Range.Delete()

# Paragraph 38: ".otherwise(\"Similar\")" -> new .when(...rlike...) line,
# followed by a brand new ".otherwise(\"Unknown\")" paragraph.
Replace-ParaText 38 '    .otherwise("Similar")' '    .when(col("cleaned_name_modified").rlike("^[A-Z]+$"), "Dissimilar")  # If word is fully randomized'
Insert-ParaAfter 38 '    .otherwise("Unknown")'

# Paragraph 37: ".when(col(\"Destination\").contains(\"RANDOM_\"), \"Dissimilar\")"
Replace-ParaText 37 `
    '.when(col("Destination").contains("RANDOM_"), "Dissimilar")' `
    '.when(col("cleaned_name_modified") != col("cleaned_name"), "Similar")'

# Paragraph 36: "when(col(\"Destination\") == col(\"Source\"), \"Same\")"
Replace-ParaText 36 `
    'when(col("Destination") == col("Source"), "Same")' `
    'when(col("cleaned_name_modified") == col("cleaned_name"), "Same")'

# Paragraph 31: ".otherwise(introduce_typos_udf(col(\"Destination\"), \"dissimilar\"))"
Replace-ParaText 31 `
    '.otherwise(introduce_typos_udf(col("Destination"), "dissimilar"))' `
    '.otherwise(introduce_typos_udf(col("cleaned_name"), "dissimilar"))                        # Completely change'

# Paragraph 30: ".when((col(\"row_number\") % 3 == 1), introduce_typos_udf(col(\"Destination\"), \"similar\"))"
Replace-ParaText 30 `
    '.when((col("row_number") % 3 == 1), introduce_typos_udf(col("Destination"), "similar"))' `
    '.when((col("row_number") % 3 == 1), introduce_typos_udf(col("cleaned_name"), "similar"))  # Minor typo'

# Paragraph 29: "when((col(\"row_number\") % 3 == 0), introduce_typos_udf(col(\"Destination\"), \"same\"))"
Replace-ParaText 29 `
    'when((col("row_number") % 3 == 0), introduce_typos_udf(col("Destination"), "same"))' `
    'when((col("row_number") % 3 == 0), introduce_typos_udf(col("cleaned_name"), "same"))      # No change'

# Paragraph 28: "    \"Destination\","
Replace-ParaText 28 '"Destination",' '"cleaned_name_modified",'

# Paragraph 26: "# Introduce typos based on random selection"
Replace-ParaText 26 '# Introduce typos based on random selection' '# Apply transformations to create modified names'

# Paragraph 24: window_spec line
Replace-ParaText 24 `
    'window_spec = Window.orderBy("Source")  # Shuffle records' `
    'window_spec = Window.orderBy("cleaned_name")  # Ensures fair distribution'

# Paragraph 23: "# Assign a row number to shuffle records equally"
Replace-ParaText 23 '# Assign a row number to shuffle records equally' '# Assign a row number to distribute changes equally'

# Paragraph 18: dissimilar-branch return line -> comment, then new return line inserted after
Replace-ParaText 18 `
    'return "RANDOM_" + str(random.randint(100, 999))  # Completely change word' `
    '# Completely change the word but keep the same length'
Insert-ParaAfter 18 '        return "".join(random.choices(string.ascii_uppercase, k=len(word)))'

# Paragraph 14: remove the trailing comment from the idx= line
Replace-ParaText 14 `
    '(word) - 2)  # Pick a random index to swap' `
    '(word) - 2)  '

# New paragraph inserted right after paragraph 13 ("elif change_type == ""similar"":")
Insert-ParaAfter 13 '        # Introduce minor typo (swap two adjacent letters)'

# Paragraph 7: docstring
Replace-ParaText 7 `
    '"""Randomly introduce typos or completely change the word."""' `
    '"""Introduce minor typos for ''Similar'' or completely change for ''Dissimilar''."""'

# Paragraph 5: "# Function to introduce typos" -> "# Initialize Spark session"
Replace-ParaText 5 '# Function to introduce typos' '# Initialize Spark session'

# New paragraphs inserted after paragraph 5, in order
Insert-ParaAfter 5 'spark = SparkSession.builder.appName("DataModification").getOrCreate()'
Insert-ParaAfter 6 '# Sample Data'
Insert-ParaAfter 7 'data = [("OTHER",), ("SAMPLE",), ("TEST",), ("DATA",), ("EXAMPLE",), ("HELLO",), ("WORLD",)]'
Insert-ParaAfter 8 'df = spark.createDataFrame(data, ["cleaned_name"])'
Insert-ParaAfter 9 '# Function to introduce minor typos or completely change a word'

# New paragraphs inserted after paragraph 1 ("import random")
Insert-ParaAfter 1 'import string'
Insert-ParaAfter 2 'from pyspark.sql import SparkSession'

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
